# blueprint_0324 -> 0325 "replay jonghyuk behavior experiment"
#
# Two independent edits:
#   1. Every "Date Placeholder" (PlaceholderFormat.Type = 16, the cached
#      datetimeFigureOut field) on the slide master and on each of its
#      custom layouts rolled from "2021. 3. 24." to "2021. 3. 25.".
#   2. On slide 1, the "Sequence1" callout textbox's second line changed
#      from "<bullet> 4<seconds>" to "<bullet> 6<seconds>".

$p = $ppt.ActivePresentation
$bullet = [string][char]0x2022

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Length -eq 12) {
                    $tr.Text = "2021. 3. 25."
                }
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# Every layout inherits/duplicates its own copy of the placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

# Slide 1 ("TextBox 85", shape id 86): "<bullet> 4<seconds>" -> "<bullet> 6<seconds>"
# Target only the "<bullet> 4" run (3 characters, starting right after the
# "Sequence1(Se1) (4-1-3-2-4)" paragraph) so the trailing "<seconds>" run is
# left untouched.
$slide1 = $p.Slides.Item(1)
$seqBox = $slide1.Shapes.Item(6)
$tr = $seqBox.TextFrame.TextRange
$run = $tr.Characters(30, 3)
if ($run.Text.Substring(2, 1) -eq "4") {
    $run.Text = "$bullet 6"
}
